# Weekly update: a new price-report entry is added for
# "Comercializadora del Agro de Limarí - Arveja Verde".
#
# The new record re-uses the Volumen/Precio data that was in row 13
# (200 / 28000 / 30000 / 29000 / 1160) but is reported for a later date
# (2022-06-30). We duplicate row 13 (Copy + Insert, which shifts every
# following row down by one and preserves formatting/styles), then
# correct the date on the freshly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Copy()
$ws.Rows.Item(13).Insert()

$ws.Range("D13").Value = "2022-06-30"
